$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 443 (pushes existing rows 443:522 down to 444:523,
# carrying their formatting with them, and extends the used range/dimension
# to A1:R523 automatically).
$ws.Rows.Item(443).Insert()

# Populate the newly inserted row 443 with the new data record.
$ws.Cells.Item(443, 1).Value = 6
$ws.Cells.Item(443, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(443, 3).Value = "Metropolitana"
$ws.Cells.Item(443, 4).Value = 44951
$ws.Cells.Item(443, 5).Value = 13
$ws.Cells.Item(443, 6).Value = 100112032
$ws.Cells.Item(443, 7).Value = "Zapallo italiano"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Primera"
$ws.Cells.Item(443, 10).Value = 1000
$ws.Cells.Item(443, 11).Value = 6000
$ws.Cells.Item(443, 12).Value = 7000
$ws.Cells.Item(443, 13).Value = 6450
$ws.Cells.Item(443, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(443, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(443, 16).Value = 129
$ws.Cells.Item(443, 17).Value = 50
$ws.Cells.Item(443, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date-formatted style as the
# cell above it (the row insert should already have copied this down, but
# set it explicitly in case a blank style slipped in).
$ws.Cells.Item(443, 4).NumberFormat = $ws.Cells.Item(442, 4).NumberFormat
